$d = $word.ActiveDocument

# Locate the paragraph that ends with the blue "You might want to replace m: else by m:else" hint.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*You might want to replace m: else by m:else*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Target paragraph (extraSpaceInElse hint) not found"
}

# Build a range that stops right before the paragraph mark, then collapse to its end
# so new runs are appended at the very end of the paragraph (same story as the rest).
$r = $target.Range.Duplicate
$r.MoveEnd(1, -1)
$r.Collapse(0)

# --- Run 1: four plain spaces (no special formatting) ---
$start = $r.Start
$r.InsertAfter("    ")
$r.Collapse(0)

# --- Run 2: "<---" in red, 16pt, light-gray highlight ---
$start = $r.Start
$r.InsertAfter("<---")
$run2 = $d.Range($start, $start + 4)
$run2.Font.Color = 255
$run2.Font.Size = 16
$run2.Font.HighlightColorIndex = 16
$r.Collapse(0)

# --- Run 3: "missing expression" in red, 16pt, light-gray highlight ---
$start = $r.Start
$r.InsertAfter("missing expression")
$run3 = $d.Range($start, $start + 19)
$run3.Font.Color = 255
$run3.Font.Size = 16
$run3.Font.HighlightColorIndex = 16
